$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text updated for both rows ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2:F3").Value = "Handed back: in sync with en-US"
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2:C3").Value = "Handed back: in sync with en-US"

# Clear hyperlinks so we can rebuild them in the exact insertion order Excel would use
$zh.Cells.Hyperlinks.Delete()

$zh.Range("I2").Value = "badf8fad-6014-4fc5-b707-bbe51ff1d588.md"
$zh.Range("J2").Value = "badf8fad-6014-4fc5-b707-bbe51ff1d588.25d770643a8d47ebd03f2530f7e96670d40d2e5c.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-18 13:03:52"

$zh.Range("I3").Value = "badf8fad-6014-4fc5-b707-bbe51ff1d588.md"
$zh.Range("J3").Value = "badf8fad-6014-4fc5-b707-bbe51ff1d588.25d770643a8d47ebd03f2530f7e96670d40d2e5c.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-18 13:03:52"

# Rebuild hyperlinks in row order (A2, I2, A3, I3) so relationship ids come out in that order
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a840dac725438d477fa4d592aee5df6ad621e4a/e2e/badf8fad-6014-4fc5-b707-bbe51ff1d588.md", "", "", "badf8fad-6014-4fc5-b707-bbe51ff1d588.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a840dac725438d477fa4d592aee5df6ad621e4a/e2e/badf8fad-6014-4fc5-b707-bbe51ff1d588.md", "", "", "badf8fad-6014-4fc5-b707-bbe51ff1d588.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a840dac725438d477fa4d592aee5df6ad621e4a/e2e/ffff5f715edf-5810-44d7-b0fa-87b70c41835e.md", "", "", "ffff5f715edf-5810-44d7-b0fa-87b70c41835e.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a840dac725438d477fa4d592aee5df6ad621e4a/e2e/badf8fad-6014-4fc5-b707-bbe51ff1d588.md", "", "", "badf8fad-6014-4fc5-b707-bbe51ff1d588.md")

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2:C3").Value = "Handed back: in sync with en-US"

$de.Cells.Hyperlinks.Delete()

$de.Range("I2").Value = "badf8fad-6014-4fc5-b707-bbe51ff1d588.md"
$de.Range("J2").Value = "badf8fad-6014-4fc5-b707-bbe51ff1d588.25d770643a8d47ebd03f2530f7e96670d40d2e5c.de-de.xlf"
$de.Range("K2").Value = "2016-08-18 13:04:16"

$de.Range("I3").Value = "badf8fad-6014-4fc5-b707-bbe51ff1d588.md"
$de.Range("J3").Value = "badf8fad-6014-4fc5-b707-bbe51ff1d588.25d770643a8d47ebd03f2530f7e96670d40d2e5c.de-de.xlf"
$de.Range("K3").Value = "2016-08-18 13:04:16"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a840dac725438d477fa4d592aee5df6ad621e4a/e2e/badf8fad-6014-4fc5-b707-bbe51ff1d588.md", "", "", "badf8fad-6014-4fc5-b707-bbe51ff1d588.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a840dac725438d477fa4d592aee5df6ad621e4a/e2e/badf8fad-6014-4fc5-b707-bbe51ff1d588.md", "", "", "badf8fad-6014-4fc5-b707-bbe51ff1d588.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a840dac725438d477fa4d592aee5df6ad621e4a/e2e/ffff5f715edf-5810-44d7-b0fa-87b70c41835e.md", "", "", "ffff5f715edf-5810-44d7-b0fa-87b70c41835e.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a840dac725438d477fa4d592aee5df6ad621e4a/e2e/badf8fad-6014-4fc5-b707-bbe51ff1d588.md", "", "", "badf8fad-6014-4fc5-b707-bbe51ff1d588.md")

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
